$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E12: {'float', 'any', 'int'} -> {'int', 'any', 'float'}
$ws.Range("E12").Value = "{'int', 'any', 'float'}"

# E13: float -> int
$ws.Range("E13").Value = "int"

# Row 69: "Scalpel Accuracy:" label and its value move from C/D to E/F
$ws.Range("C69").Value = $null
$ws.Range("D69").Value = $null
$ws.Range("E69").Value = "Scalpel Accuracy:"
$ws.Range("F69").Value = 100

# Row 70: fix label text
$ws.Range("E70").Value = "Accuracy vs PyType"
